# Update "想去人数" (want-to-go count) values in column F across sheets,
# matching the output generated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 1658
$ws1.Range("F3").Value  = 865
$ws1.Range("F7").Value  = 809
$ws1.Range("F9").Value  = 1534
$ws1.Range("F10").Value = 308
$ws1.Range("F13").Value = 77
$ws1.Range("F16").Value = 515
$ws1.Range("F17").Value = 65
$ws1.Range("F22").Value = 581
$ws1.Range("F26").Value = 783

# Sheet "演出" (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 1038
$ws2.Range("F7").Value = 153

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 1658
$ws4.Range("F5").Value  = 865
$ws4.Range("F7").Value  = 1038
$ws4.Range("F10").Value = 809
$ws4.Range("F12").Value = 1534
$ws4.Range("F13").Value = 308
$ws4.Range("F16").Value = 77
$ws4.Range("F19").Value = 515
$ws4.Range("F20").Value = 65
$ws4.Range("F28").Value = 153
$ws4.Range("F29").Value = 153
$ws4.Range("F30").Value = 581
$ws4.Range("F34").Value = 783
